$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns with refreshed crypto
# quotes. Several Price values (e.g. "1.010", "344.04") look like plain
# numbers to Excel, so a leading apostrophe is used to force them to stay
# as text (matching the source data, which stores these as strings) instead
# of being silently reinterpreted as numeric values.
$ws.Range("D2").Value = '29.727.57'
$ws.Range("E2").Value = '  -2.60%  '
$ws.Range("D3").Value = '2.095.84'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").Value = '''1.010'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''344.04'
$ws.Range("E5").Value = '  -2.40%  '
$ws.Range("D6").Value = '''1.009'
$ws.Range("D7").Value = '''0.5174'
$ws.Range("E7").Value = '  -1.62%  '
$ws.Range("D8").Value = '''0.4391'
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("D9").Value = '''53.04'
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").Value = '''0.09262'
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").Value = '''1.165'
$ws.Range("E11").Value = '  -2.43%  '
$ws.Range("D12").Value = '''24.99'
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").Value = '2.101.02'
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("D14").Value = '''8.278'
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Value = '''6.775'
$ws.Range("E15").Value = '  -1.72%  '
$ws.Range("D16").Value = '''99.63'
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("D17").Value = '''0.00001152'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").Value = '''1.010'
$ws.Range("D19").Value = '''20.84'
$ws.Range("E19").Value = '  +1.57%  '
$ws.Range("D20").Value = '''0.06641'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D22").Value = '''6.207'
$ws.Range("E22").Value = '  -2.63%  '
$ws.Range("D23").Value = '29.766.86'
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("D24").Value = '''12.50'
$ws.Range("E24").Value = '  -2.78%  '
$ws.Range("D25").Value = '''2.317'
$ws.Range("E25").Value = '  -2.63%  '
$ws.Range("D26").Value = '2.348.84'
$ws.Range("E26").Value = '  -1.81%  '
$ws.Range("D27").Value = '''21.94'
$ws.Range("E27").Value = '  -2.80%  '
$ws.Range("D28").Value = '''2.521'
$ws.Range("E28").Value = '  -3.20%  '
$ws.Range("D29").Value = '''161.46'
$ws.Range("E29").Value = '  -2.22%  '
$ws.Range("D30").Value = '''133.13'
$ws.Range("E30").Value = '  -1.97%  '
$ws.Range("D31").Value = '''1.141'
$ws.Range("E31").Value = '  -6.85%  '
$ws.Range("E32").Value = '  -3.11%  '
$ws.Range("D33").Value = '''1.645'
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("D34").Value = '''6.182'
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("D35").Value = '''3.955'
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("D36").Value = '''6.316'
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("D37").Value = '''10.23'
$ws.Range("E37").Value = '  -1.88%  '
$ws.Range("D38").Value = '''0.02579'
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("D39").Value = '''0.06739'
$ws.Range("E39").Value = '  -3.50%  '
$ws.Range("D40").Value = '''0.7001'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("E41").Value = '  -1.89%  '
$ws.Range("D42").Value = '''0.2229'
$ws.Range("E42").Value = '  -5.41%  '
$ws.Range("D43").Value = '''1.323'
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("D44").Value = '''0.6925'
$ws.Range("E44").Value = '  +6.34%  '
$ws.Range("D45").Value = '''14.32'
$ws.Range("E45").Value = '  -3.00%  '
$ws.Range("D46").Value = '''2.323'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("E48").Value = '  -5.35%  '
$ws.Range("D49").Value = '''1.219'
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("D50").Value = '''82.11'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("E51").Value = '  -2.41%  '
